$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New monthly_mean (column D) and weekly_share (column F) values for rows 2..54,
# recomputed after the monthly_mean denominator was corrected (weekly_share = weekly_mean / monthly_mean).
$monthlyMean = @(101.285, 101.285, 101.285, 101.285, 101.285, 95.556, 95.556, 95.556, 95.556, 95.556, 114.613, 114.613, 114.613, 114.613, 117.204, 117.204, 117.204, 117.204, 81.501, 81.501, 81.501, 81.501, 81.501, 75.141, 75.141, 75.141, 75.141, 29.892, 29.892, 29.892, 29.892, 30.37, 30.37, 30.37, 30.37, 30.37, 44.683, 44.683, 44.683, 44.683, 111.232, 111.232, 111.232, 111.232, 111.232, 127.063, 127.063, 127.063, 127.063, 145.175, 145.175, 145.175, 145.175)

$weeklyShare = @(0.796, 1.071, 1.164, 1.218, 0.751, 0.654, 1.065, 1.176, 0.876, 1.23, 1.055, 1.148, 1.009, 0.788, 0.946, 1.167, 1.101, 0.786, 0.987, 1.168, 1.044, 0.885, 0.916, 1.425, 1.063, 0.763, 0.749, 1.3, 1.017, 0.911, 0.772, 0.917, 1.174, 1.26, 0.877, 0.772, 0.467, 0.621, 1.096, 1.816, 0.481, 0.918, 0.809, 1.766, 1.026, 0.949, 1.032, 0.819, 1.2, 1.343, 1.095, 0.759, 0.804)

$firstRow = 2
for ($i = 0; $i -lt $monthlyMean.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 4).Value = $monthlyMean[$i]
    $ws.Cells.Item($row, 6).Value = $weeklyShare[$i]
}
